# The deck currently uses the "Integral" theme (ppt/theme/theme2.xml) for
# its slide master / slides, while ppt/theme/theme1.xml ("Office Theme")
# is only wired to the notes master. The edit being replayed here swaps the
# two theme parts: the slides revert to the stock "Office Theme" palette
# and the (previously unused) "Integral" palette is swapped in on the other
# slot. We reproduce the user-visible effect - the deck's active theme
# colors changing from the Integral palette to the Office palette - via the
# supported PowerPoint theme-color COM surface.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeColor($tcs, $index, $r, $g, $b) {
    $rgb = $r + ($g * 256) + ($b * 65536)
    $tcs.Colors($index).RGB = $rgb
}

# Target palette: the default Office theme colors (what ppt/theme/theme2.xml
# should contain after the swap), applied in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
Set-ThemeColor $tcs 1  0   0   0     # dk1       000000
Set-ThemeColor $tcs 2  255 255 255   # lt1       FFFFFF
Set-ThemeColor $tcs 3  68  84  106   # dk2       44546A
Set-ThemeColor $tcs 4  231 230 230   # lt2       E7E6E6
Set-ThemeColor $tcs 5  91  155 213   # accent1   5B9BD5
Set-ThemeColor $tcs 6  237 125 49    # accent2   ED7D31
Set-ThemeColor $tcs 7  165 165 165   # accent3   A5A5A5
Set-ThemeColor $tcs 8  255 192 0     # accent4   FFC000
Set-ThemeColor $tcs 9  68  114 196   # accent5   4472C4
Set-ThemeColor $tcs 10 112 173 71    # accent6   70AD47
Set-ThemeColor $tcs 11 5   99  193   # hlink     0563C1
Set-ThemeColor $tcs 12 149 79  114   # folHlink  954F72
